$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.394.16'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.891.59'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.45'
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('E7').Value = '  -0.82%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '43.33'
$ws.Range('E8').Value = '  +5.17%  '
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '53.82'
$ws.Range('E10').Value = '  +1.73%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0744'
$ws.Range('E11').Value = '  -1.83%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0969'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '13.31'
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('D14').Value = '2.165.29'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.760'
$ws.Range('E15').Value = '  +3.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.90'
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').Value = '1.884.13'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = '35.468.92'
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '73.23'
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '245.03'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('E23').Value = '  -1.92%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.68'
$ws.Range('E24').Value = '  +10.42%  '
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -7.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '166.17'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('E29').Value = '  -0.74%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.128'
$ws.Range('E30').Value = '  -1.44%  '
$ws.Range('D31').Value = '4.128.47'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.74'
$ws.Range('E32').Value = '  +10.83%  '
$ws.Range('E33').Value = '  -1.27%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0585'
$ws.Range('E34').Value = '  -3.54%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.17'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.86'
$ws.Range('E36').Value = '  -12.78%  '
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.848'
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('E40').Value = '  +7.19%  '
$ws.Range('E41').Value = '  +2.78%  '
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '96.95'
$ws.Range('E43').Value = '  -3.06%  '
$ws.Range('E44').Value = '  -2.70%  '
$ws.Range('D45').Value = '1.300.27'
$ws.Range('E45').Value = '  -2.26%  '
$ws.Range('E46').Value = '  -4.39%  '
$ws.Range('E47').Value = '  +7.42%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.41'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('E49').Value = '  +3.93%  '
$ws.Range('E50').Value = '  -0.46%  '
$ws.Range('E51').Value = '  -5.28%  '
